# New Submission Synced: 2026-02-08 19:27:05
# Target sheet: "JSS 3E" (18th tab) of the results workbook.
#   - Row 6, column C ("Admission No") was stored as text "4"; it is
#     corrected to the real number 4.
#   - A brand-new submission row (row 7) is appended with the four
#     form columns: Timestamp, Full Name, Admission No, AI Score.
#     The Admission No "33" is synced as literal text (matching how the
#     prior rows' raw submissions arrive), while the AI Score is numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# --- Fix existing row 6: Admission No should be a real number, not text ---
$ws.Range("C6").Value = 4

# --- Append the newly synced submission as row 7 ---
$ws.Range("A7").Value = "2026-02-08 19:27:04"
$ws.Range("B7").Value = "FATIMA BUKAR WAZIRI"

# Admission No "33" must stay text (leading apostrophe forces text entry,
# then the number format / style is reset to plain/default afterward so no
# stray "Text" formatting is left behind on the cell).
$ws.Range("C7").Value = "'33"
$ws.Range("C7").NumberFormat = "General"
$ws.Range("C7").Style = "Normal"

$ws.Range("D7").Value = 8
